# Update the "Metadata" worksheet (first sheet) to reflect the new
# IG version/publish metadata:
#   - Version bumped from 5.0.0 to 6.0.0
#   - Date bumped to the new publish timestamp
#   - Publisher value filled in ("Alvearie Team")
#   - The old duplicated "Contact" / "No display for ContactDetail" rows
#     are replaced by a single "Jurisdiction" / "United States of America" row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Version value (row 3)
$ws.Range("B3").Value = "6.0.0"

# Date value (row 8)
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (row 9) - was blank
$ws.Range("B9").Value = "Alvearie Team"

# Turn the first of the two duplicate "Contact" rows (row 10) into the
# new "Jurisdiction" row
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Remove the now-redundant duplicate "Contact" row (row 11), shifting
# everything below it up by one row
$ws.Rows.Item(11).Delete()
